$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "'67.392.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "  +0.53%  "

$ws.Range("D3").Value2 = "'2.476.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "  +0.19%  "

$ws.Range("D4").Value2 = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  -0.04%  "

$ws.Range("D5").Value2 = "'584.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +0.14%  "

$ws.Range("D6").Value2 = "'175.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +3.68%  "

$ws.Range("E7").Value2 = "  +0.02%  "

$ws.Range("E8").Value2 = "  -0.29%  "

$ws.Range("D9").Value2 = "'0.139"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +2.71%  "

$ws.Range("E10").Value2 = "  +0.75%  "

$ws.Range("B11").Value2 = "Cardano"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value2 = "'0.336"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +1.47%  "

$ws.Range("B12").Value2 = "Toncoin"
$ws.Range("C12").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value2 = "'4.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "  -0.04%  "

$ws.Range("D13").Value2 = "'2.928.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +0.18%  "

$ws.Range("D14").Value2 = "'25.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  -0.39%  "

$ws.Range("D15").Value2 = "'67.241.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +0.50%  "

$ws.Range("D16").Value2 = "'0.0000170"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +0.28%  "

$ws.Range("D17").Value2 = "'2.479.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  -0.30%  "

$ws.Range("D18").Value2 = "'10.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  -1.29%  "

$ws.Range("D19").Value2 = "'7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  -0.93%  "

$ws.Range("D20").Value2 = "'350.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  -0.59%  "

$ws.Range("D21").Value2 = "'4.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "  -0.13%  "

$ws.Range("E22").Value2 = "  -0.01%  "

$ws.Range("D23").Value2 = "'70.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +1.82%  "

$ws.Range("D24").Value2 = "'4.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "  -0.96%  "

$ws.Range("D25").Value2 = "'1.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  -1.32%  "

$ws.Range("D26").Value2 = "'9.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +0.20%  "

$ws.Range("D27").Value2 = "'2.595.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  -0.10%  "

$ws.Range("D28").Value2 = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  -0.02%  "

$ws.Range("D29").Value2 = "'0.0₃0907"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  -0.12%  "

$ws.Range("D30").Value2 = "'503.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  -1.60%  "

$ws.Range("D31").Value2 = "'7.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +0.84%  "

$ws.Range("D32").Value2 = "'1.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +0.33%  "

$ws.Range("D33").Value2 = "'1.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  -0.71%  "

$ws.Range("E34").Value2 = "  +0.00%  "

$ws.Range("E35").Value2 = "  +3.43%  "

$ws.Range("D36").Value2 = "'162.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +1.51%  "

$ws.Range("D37").Value2 = "'18.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "  +0.21%  "

$ws.Range("D38").Value2 = "'18.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  -0.63%  "

$ws.Range("E39").Value2 = "  -0.73%  "

$ws.Range("E40").Value2 = "  +0.04%  "

$ws.Range("D41").Value2 = "'1.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +1.05%  "

$ws.Range("D42").Value2 = "'0.329"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +0.32%  "

$ws.Range("D43").Value2 = "'4.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  -0.06%  "

$ws.Range("D44").Value2 = "'2.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +0.83%  "

$ws.Range("D45").Value2 = "'143.20"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +1.81%  "

$ws.Range("D46").Value2 = "'3.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "  +1.02%  "

$ws.Range("D47").Value2 = "'0.0₆0259"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +1.12%  "

$ws.Range("D48").Value2 = "'0.510"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  -1.03%  "

$ws.Range("D49").Value2 = "'0.0743"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +1.75%  "

$ws.Range("D50").Value2 = "'1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  -0.95%  "

$ws.Range("D51").Value2 = "'0.583"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  -0.04%  "
